# Weekly fruit/vegetable price update: a new weekly price observation
# was inserted into the "Espinaca" (Mercado Mayorista Lo Valledor de
# Santiago) dataset. This pushes all subsequent rows down by one and
# grows the used range from A1:R552 to A1:R553.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 516, shifting rows 516:552 to 517:553.
$ws.Rows.Item(516).Insert()

# Populate the newly inserted row 516 with the new weekly observation.
$ws.Range("A516").Value = 6
$ws.Range("B516").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C516").Value = "Metropolitana"
$ws.Range("D516").Value = 44714
$ws.Range("E516").Value = 13
$ws.Range("F516").Value = 100112012
$ws.Range("G516").Value = "Espinaca"
$ws.Range("H516").Value = "Sin especificar"
$ws.Range("I516").Value = "Primera"
$ws.Range("J516").Value = 450
$ws.Range("K516").Value = 4000
$ws.Range("L516").Value = 4500
$ws.Range("M516").Value = 4211
$ws.Range("N516").Value = "`$/cuna 10 kilos"
$ws.Range("O516").Value = "Región Metropolitana"
$ws.Range("P516").Value = 421
$ws.Range("Q516").Value = 10
$ws.Range("R516").Value = "Hortaliza"
